{"js": "// Commit: \"Equality, I put ~= for not equal\"\n// The \"Equality ==\" bullet line gets \", ~=\" appended, becoming\n// \"Equality ==, ~=\" (documenting the not-equal operator alongside ==).\n\nconst results = context.document.body.search(\"Equality ==\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the 'Equality ==' line in the document\");\n}\n\n// There should be exactly one such line (the \"Operators\" bullet list).\nconst target = results.items[0];\ntarget.insertText(\", ~=\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Commit: \"Equality, I put ~= for not equal\"\n#\n# The \"Operators:\" bullet list has a line reading \"Equality ==\".\n# Append \", ~=\" to it so it documents the not-equal operator too,\n# producing \"Equality ==, ~=\".\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Equality ==\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n\nif ($find.Execute()) {\n    # $rng now spans the found text (\"Equality ==\"); collapse to its\n    # end point (0 = wdCollapseEnd) and insert the addition right after.\n    $rng.Collapse(0)\n    $rng.InsertAfter(\", ~=\")\n} else {\n    throw \"Could not find the 'Equality ==' line in the document\"\n}\n"}
